# Updates the "cryptos" price table (rows 2-51) to reflect the refreshed
# coinranking.com snapshot described in the commit message. Each row is a
# coin; columns are B=Coin, C=Link, D=Price, E=Volume(1h). Most rows only get
# new Price/Volume figures, but a few coin pairs also swapped row order
# (rows 32/33, 41/42, 47/48), which is expressed here as updates to B/C/D/E
# for both rows involved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.730.74"
$ws.Range("E2").Value = "  +0.49%  "
# Row 3
$ws.Range("D3").Value = "1.701.36"
$ws.Range("E3").Value = "  +0.28%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.29%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.47"
$ws.Range("E5").Value = "  -0.65%  "
# Row 6
$ws.Range("E6").Value = "  +0.22%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3943"
$ws.Range("E7").Value = "  -0.32%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4044"
$ws.Range("E8").Value = "  +0.60%  "
# Row 9
$ws.Range("E9").Value = "  -0.16%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.005"
$ws.Range("E10").Value = "  +0.24%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.62"
$ws.Range("E11").Value = "  +0.53%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08870"
$ws.Range("E12").Value = "  +0.90%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.542"
$ws.Range("E13").Value = "  +3.73%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.69"
$ws.Range("E14").Value = "  +1.71%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.185"
$ws.Range("E15").Value = "  +7.20%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001326"
$ws.Range("E16").Value = "  +0.32%  "
# Row 17
$ws.Range("D17").Value = "1.707.35"
$ws.Range("E17").Value = "  +0.61%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.80"
$ws.Range("E18").Value = "  -1.30%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07047"
$ws.Range("E19").Value = "  +0.55%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.74"
$ws.Range("E20").Value = "  +0.06%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.078"
$ws.Range("E21").Value = "  +2.55%  "
# Row 22
$ws.Range("E22").Value = "  +0.09%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.66"
$ws.Range("E23").Value = "  +3.85%  "
# Row 24
$ws.Range("D24").Value = "24.716.32"
$ws.Range("E24").Value = "  +0.43%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.184"
$ws.Range("E25").Value = "  +3.67%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369"
$ws.Range("E26").Value = "  +1.29%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("E27").Value = "  +1.68%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.28"
$ws.Range("E28").Value = "  +1.62%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.715"
$ws.Range("E29").Value = "  +16.69%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.86"
$ws.Range("E30").Value = "  +0.94%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.181"
$ws.Range("E31").Value = "  -1.20%  "
# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09047"
$ws.Range("E32").Value = "  +5.96%  "
# Row 33
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.678"
$ws.Range("E33").Value = "  +2.05%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.076"
$ws.Range("E34").Value = "  -2.73%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.993"
$ws.Range("E35").Value = "  +0.98%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.12"
$ws.Range("E36").Value = "  -3.18%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2759"
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.54"
$ws.Range("E38").Value = "  -0.47%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02780"
$ws.Range("E39").Value = "  +0.14%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09173"
$ws.Range("E40").Value = "  +1.46%  "
# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.465"
$ws.Range("E41").Value = "  -0.10%  "
# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7705"
$ws.Range("E42").Value = "  -0.10%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.05"
$ws.Range("E43").Value = "  +3.47%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7198"
$ws.Range("E44").Value = "  -0.49%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.576"
$ws.Range("E45").Value = "  +1.94%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.220"
$ws.Range("E46").Value = "  -0.38%  "
# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.23%  "
# Row 48
$ws.Range("B48").Value = "Flow"
$ws.Range("C48").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.348"
$ws.Range("E48").Value = "  -0.50%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.71"
$ws.Range("E49").Value = "  -0.36%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.24"
$ws.Range("E50").Value = "  +3.18%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07985"
$ws.Range("E51").Value = "  -0.75%  "
